$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): insert new column headers for G1:Q1 ---
# Columns G:L already existed (style s="1" carries over automatically).
# Columns M:Q are brand new (outside the old A1:L31 dimension) and need
# the header style (bold, centered, bordered) copied over explicitly.
$ws.Cells.Item(1, 7).Value = "vet_costs"
$ws.Cells.Item(1, 8).Value = "bedding_litter_costs"
$ws.Cells.Item(1, 9).Value = "marketing_costs"
$ws.Cells.Item(1, 10).Value = "service_costs"
$ws.Cells.Item(1, 11).Value = "utility_costs"
$ws.Cells.Item(1, 12).Value = "repair_costs"
$ws.Cells.Item(1, 13).Value = "total_operating_costs"
$ws.Cells.Item(1, 14).Value = "overhead"
$ws.Cells.Item(1, 15).Value = "total_costs"
$ws.Cells.Item(1, 16).Value = "cows_per_farm"
$ws.Cells.Item(1, 17).Value = "output_per_cow"

# Copy the header formatting (bold/border/centered) from an existing
# header cell onto the newly-created M1:Q1 header cells.
$ws.Range("A1").Copy()
$ws.Range("M1:Q1").PasteSpecial(-4122)

# --- Data rows 2-31: populate the new cost-breakdown columns G:Q ---
# Row 2
$ws.Cells.Item(2, 7).Value = 0.54
$ws.Cells.Item(2, 8).Value = 0.38
$ws.Cells.Item(2, 9).Value = 0.43
$ws.Cells.Item(2, 10).Value = 0.49
$ws.Cells.Item(2, 11).Value = 0.89
$ws.Cells.Item(2, 12).Value = 1.07
$ws.Cells.Item(2, 13).Value = 14.3
$ws.Cells.Item(2, 14).Value = 16.31
$ws.Cells.Item(2, 15).Value = 30.61
$ws.Cells.Item(2, 16).Value = 59
$ws.Cells.Item(2, 17).Value = 13058

# Row 3
$ws.Cells.Item(3, 7).Value = 0.35
$ws.Cells.Item(3, 8).Value = 0.37
$ws.Cells.Item(3, 9).Value = 0.24
$ws.Cells.Item(3, 10).Value = 0.1
$ws.Cells.Item(3, 11).Value = 1.12
$ws.Cells.Item(3, 12).Value = 0.9399999999999999
$ws.Cells.Item(3, 13).Value = 16.63
$ws.Cells.Item(3, 14).Value = 21.37
$ws.Cells.Item(3, 15).Value = 38
$ws.Cells.Item(3, 16).Value = 53
$ws.Cells.Item(3, 17).Value = 10956

# Row 4
$ws.Cells.Item(4, 7).Value = 0.68
$ws.Cells.Item(4, 8).Value = 0.46
$ws.Cells.Item(4, 9).Value = 0.31
$ws.Cells.Item(4, 10).Value = 0.23
$ws.Cells.Item(4, 11).Value = 0.73
$ws.Cells.Item(4, 12).Value = 0.9399999999999999
$ws.Cells.Item(4, 13).Value = 17.1
$ws.Cells.Item(4, 14).Value = 14.61
$ws.Cells.Item(4, 15).Value = 31.71
$ws.Cells.Item(4, 16).Value = 46
$ws.Cells.Item(4, 17).Value = 12177

# Row 5
$ws.Cells.Item(5, 7).Value = 0.8100000000000001
$ws.Cells.Item(5, 8).Value = 0.52
$ws.Cells.Item(5, 9).Value = 0.24
$ws.Cells.Item(5, 10).Value = 0.22
$ws.Cells.Item(5, 11).Value = 0.98
$ws.Cells.Item(5, 12).Value = 0.87
$ws.Cells.Item(5, 13).Value = 17.62
$ws.Cells.Item(5, 14).Value = 18.88
$ws.Cells.Item(5, 15).Value = 36.5
$ws.Cells.Item(5, 16).Value = 62
$ws.Cells.Item(5, 17).Value = 11982

# Row 6
$ws.Cells.Item(6, 7).Value = 0.6
$ws.Cells.Item(6, 8).Value = 0.39
$ws.Cells.Item(6, 9).Value = 0.22
$ws.Cells.Item(6, 10).Value = 0.33
$ws.Cells.Item(6, 11).Value = 0.84
$ws.Cells.Item(6, 12).Value = 0.85
$ws.Cells.Item(6, 13).Value = 17.11
$ws.Cells.Item(6, 14).Value = 14.95
$ws.Cells.Item(6, 15).Value = 32.07
$ws.Cells.Item(6, 16).Value = 65
$ws.Cells.Item(6, 17).Value = 13101

# Row 7
$ws.Cells.Item(7, 7).Value = 0.49
$ws.Cells.Item(7, 8).Value = 0.47
$ws.Cells.Item(7, 9).Value = 0.25
$ws.Cells.Item(7, 10).Value = 0.6
$ws.Cells.Item(7, 11).Value = 1.2
$ws.Cells.Item(7, 12).Value = 0.65
$ws.Cells.Item(7, 13).Value = 17.89
$ws.Cells.Item(7, 14).Value = 21.99
$ws.Cells.Item(7, 15).Value = 39.88
$ws.Cells.Item(7, 16).Value = 36
$ws.Cells.Item(7, 17).Value = 10542

# Row 8
$ws.Cells.Item(8, 7).Value = 0.62
$ws.Cells.Item(8, 8).Value = 0.52
$ws.Cells.Item(8, 9).Value = 0.39
$ws.Cells.Item(8, 10).Value = 0.68
$ws.Cells.Item(8, 11).Value = 1.24
$ws.Cells.Item(8, 12).Value = 1.19
$ws.Cells.Item(8, 13).Value = 16.39
$ws.Cells.Item(8, 14).Value = 18.1
$ws.Cells.Item(8, 15).Value = 34.48999999999999
$ws.Cells.Item(8, 16).Value = 73
$ws.Cells.Item(8, 17).Value = 12880

# Row 9
$ws.Cells.Item(9, 7).Value = 0.51
$ws.Cells.Item(9, 8).Value = 0.44
$ws.Cells.Item(9, 9).Value = 0.23
$ws.Cells.Item(9, 10).Value = 0.41
$ws.Cells.Item(9, 11).Value = 1.44
$ws.Cells.Item(9, 12).Value = 1.55
$ws.Cells.Item(9, 13).Value = 18.81
$ws.Cells.Item(9, 14).Value = 20.52
$ws.Cells.Item(9, 15).Value = 39.33
$ws.Cells.Item(9, 16).Value = 72
$ws.Cells.Item(9, 17).Value = 12247

# Row 10
$ws.Cells.Item(10, 7).Value = 0.57
$ws.Cells.Item(10, 8).Value = 0.5
$ws.Cells.Item(10, 9).Value = 0.23
$ws.Cells.Item(10, 10).Value = 0.45
$ws.Cells.Item(10, 11).Value = 1.32
$ws.Cells.Item(10, 12).Value = 1.96
$ws.Cells.Item(10, 13).Value = 21.05
$ws.Cells.Item(10, 14).Value = 22.45
$ws.Cells.Item(10, 15).Value = 43.5
$ws.Cells.Item(10, 16).Value = 60
$ws.Cells.Item(10, 17).Value = 12145

# Row 11
$ws.Cells.Item(11, 7).Value = 0.8
$ws.Cells.Item(11, 8).Value = 0.58
$ws.Cells.Item(11, 9).Value = 0.19
$ws.Cells.Item(11, 10).Value = 0.93
$ws.Cells.Item(11, 11).Value = 1.09
$ws.Cells.Item(11, 12).Value = 1.3
$ws.Cells.Item(11, 13).Value = 19.85
$ws.Cells.Item(11, 14).Value = 21.15
$ws.Cells.Item(11, 15).Value = 41
$ws.Cells.Item(11, 16).Value = 56
$ws.Cells.Item(11, 17).Value = 12940

# Row 12
$ws.Cells.Item(12, 7).Value = 0.55
$ws.Cells.Item(12, 8).Value = 0.41
$ws.Cells.Item(12, 9).Value = 0.2
$ws.Cells.Item(12, 10).Value = 0.85
$ws.Cells.Item(12, 11).Value = 1.29
$ws.Cells.Item(12, 12).Value = 2.13
$ws.Cells.Item(12, 13).Value = 18.56
$ws.Cells.Item(12, 14).Value = 16.47
$ws.Cells.Item(12, 15).Value = 35.03
$ws.Cells.Item(12, 16).Value = 89
$ws.Cells.Item(12, 17).Value = 13376

# Row 13
$ws.Cells.Item(13, 7).Value = 0.64
$ws.Cells.Item(13, 8).Value = 0.48
$ws.Cells.Item(13, 9).Value = 0.1
$ws.Cells.Item(13, 10).Value = 0.52
$ws.Cells.Item(13, 11).Value = 1.08
$ws.Cells.Item(13, 12).Value = 1.8
$ws.Cells.Item(13, 13).Value = 28.56000000000001
$ws.Cells.Item(13, 14).Value = 18.06
$ws.Cells.Item(13, 15).Value = 46.62
$ws.Cells.Item(13, 16).Value = 62
$ws.Cells.Item(13, 17).Value = 12223

# Row 14
$ws.Cells.Item(14, 7).Value = 0.59
$ws.Cells.Item(14, 8).Value = 0.44
$ws.Cells.Item(14, 9).Value = 0.22
$ws.Cells.Item(14, 10).Value = 0.7
$ws.Cells.Item(14, 11).Value = 0.79
$ws.Cells.Item(14, 12).Value = 1.4
$ws.Cells.Item(14, 13).Value = 17.51
$ws.Cells.Item(14, 14).Value = 20.35
$ws.Cells.Item(14, 15).Value = 37.86
$ws.Cells.Item(14, 16).Value = 37
$ws.Cells.Item(14, 17).Value = 13178

# Row 15
$ws.Cells.Item(15, 7).Value = 0.5600000000000001
$ws.Cells.Item(15, 8).Value = 0.25
$ws.Cells.Item(15, 9).Value = 0.25
$ws.Cells.Item(15, 10).Value = 0.66
$ws.Cells.Item(15, 11).Value = 0.72
$ws.Cells.Item(15, 12).Value = 1.36
$ws.Cells.Item(15, 13).Value = 20
$ws.Cells.Item(15, 14).Value = 17.58
$ws.Cells.Item(15, 15).Value = 37.58
$ws.Cells.Item(15, 16).Value = 49
$ws.Cells.Item(15, 17).Value = 14543

# Row 16
$ws.Cells.Item(16, 7).Value = 0.57
$ws.Cells.Item(16, 8).Value = 0.78
$ws.Cells.Item(16, 9).Value = 0.37
$ws.Cells.Item(16, 10).Value = 0.75
$ws.Cells.Item(16, 11).Value = 1.44
$ws.Cells.Item(16, 12).Value = 3.02
$ws.Cells.Item(16, 13).Value = 22.64
$ws.Cells.Item(16, 14).Value = 18.17
$ws.Cells.Item(16, 15).Value = 40.81
$ws.Cells.Item(16, 16).Value = 63
$ws.Cells.Item(16, 17).Value = 14211

# Row 17
$ws.Cells.Item(17, 7).Value = 0.54
$ws.Cells.Item(17, 8).Value = 0.04
$ws.Cells.Item(17, 9).Value = 0.28
$ws.Cells.Item(17, 10).Value = 0.31
$ws.Cells.Item(17, 11).Value = 0.41
$ws.Cells.Item(17, 12).Value = 0.33
$ws.Cells.Item(17, 13).Value = 10.24
$ws.Cells.Item(17, 14).Value = 3.88
$ws.Cells.Item(17, 15).Value = 14.12
$ws.Cells.Item(17, 16).Value = 988
$ws.Cells.Item(17, 17).Value = 19987

# Row 18
$ws.Cells.Item(18, 7).Value = 0.79
$ws.Cells.Item(18, 8).Value = 0.03
$ws.Cells.Item(18, 9).Value = 0.18
$ws.Cells.Item(18, 10).Value = 0.72
$ws.Cells.Item(18, 11).Value = 0.7
$ws.Cells.Item(18, 12).Value = 0.73
$ws.Cells.Item(18, 13).Value = 12.42
$ws.Cells.Item(18, 14).Value = 6.42
$ws.Cells.Item(18, 15).Value = 18.84
$ws.Cells.Item(18, 16).Value = 1013
$ws.Cells.Item(18, 17).Value = 17572

# Row 19
$ws.Cells.Item(19, 7).Value = 0.62
$ws.Cells.Item(19, 8).Value = 0.03
$ws.Cells.Item(19, 9).Value = 0.7
$ws.Cells.Item(19, 10).Value = 0.79
$ws.Cells.Item(19, 11).Value = 0.59
$ws.Cells.Item(19, 12).Value = 0.5600000000000001
$ws.Cells.Item(19, 13).Value = 11.06
$ws.Cells.Item(19, 14).Value = 7.419999999999999
$ws.Cells.Item(19, 15).Value = 18.48
$ws.Cells.Item(19, 16).Value = 273
$ws.Cells.Item(19, 17).Value = 17294

# Row 20
$ws.Cells.Item(20, 7).Value = 0.65
$ws.Cells.Item(20, 8).Value = 0.22
$ws.Cells.Item(20, 9).Value = 0.25
$ws.Cells.Item(20, 10).Value = 0.52
$ws.Cells.Item(20, 11).Value = 0.55
$ws.Cells.Item(20, 12).Value = 0.63
$ws.Cells.Item(20, 13).Value = 11.93
$ws.Cells.Item(20, 14).Value = 9.010000000000002
$ws.Cells.Item(20, 15).Value = 20.94
$ws.Cells.Item(20, 16).Value = 133
$ws.Cells.Item(20, 17).Value = 17949

# Row 21
$ws.Cells.Item(21, 7).Value = 1.11
$ws.Cells.Item(21, 8).Value = 0.32
$ws.Cells.Item(21, 9).Value = 0.29
$ws.Cells.Item(21, 10).Value = 0.43
$ws.Cells.Item(21, 11).Value = 0.59
$ws.Cells.Item(21, 12).Value = 0.72
$ws.Cells.Item(21, 13).Value = 11.23
$ws.Cells.Item(21, 14).Value = 9.57
$ws.Cells.Item(21, 15).Value = 20.8
$ws.Cells.Item(21, 16).Value = 127
$ws.Cells.Item(21, 17).Value = 18389

# Row 22
$ws.Cells.Item(22, 7).Value = 0.57
$ws.Cells.Item(22, 8).Value = 0.08
$ws.Cells.Item(22, 9).Value = 0.28
$ws.Cells.Item(22, 10).Value = 0.43
$ws.Cells.Item(22, 11).Value = 0.51
$ws.Cells.Item(22, 12).Value = 0.37
$ws.Cells.Item(22, 13).Value = 11.46
$ws.Cells.Item(22, 14).Value = 4.56
$ws.Cells.Item(22, 15).Value = 16.02
$ws.Cells.Item(22, 16).Value = 999
$ws.Cells.Item(22, 17).Value = 22569

# Row 23
$ws.Cells.Item(23, 7).Value = 0.57
$ws.Cells.Item(23, 8).Value = 0.01
$ws.Cells.Item(23, 9).Value = 0.15
$ws.Cells.Item(23, 10).Value = 0.48
$ws.Cells.Item(23, 11).Value = 0.8100000000000001
$ws.Cells.Item(23, 12).Value = 0.53
$ws.Cells.Item(23, 13).Value = 14.25
$ws.Cells.Item(23, 14).Value = 8.200000000000001
$ws.Cells.Item(23, 15).Value = 22.45
$ws.Cells.Item(23, 16).Value = 905
$ws.Cells.Item(23, 17).Value = 18071

# Row 24
$ws.Cells.Item(24, 7).Value = 0.59
$ws.Cells.Item(24, 8).Value = 0.02
$ws.Cells.Item(24, 9).Value = 0.2
$ws.Cells.Item(24, 10).Value = 0.6899999999999999
$ws.Cells.Item(24, 11).Value = 0.82
$ws.Cells.Item(24, 12).Value = 0.53
$ws.Cells.Item(24, 13).Value = 12.65
$ws.Cells.Item(24, 14).Value = 9.869999999999999
$ws.Cells.Item(24, 15).Value = 22.52
$ws.Cells.Item(24, 16).Value = 188
$ws.Cells.Item(24, 17).Value = 15518

# Row 25
$ws.Cells.Item(25, 7).Value = 0.54
$ws.Cells.Item(25, 8).Value = 0.23
$ws.Cells.Item(25, 9).Value = 0.26
$ws.Cells.Item(25, 10).Value = 0.26
$ws.Cells.Item(25, 11).Value = 0.4
$ws.Cells.Item(25, 12).Value = 0.31
$ws.Cells.Item(25, 13).Value = 10.14
$ws.Cells.Item(25, 14).Value = 3.36
$ws.Cells.Item(25, 15).Value = 13.5
$ws.Cells.Item(25, 16).Value = 1110
$ws.Cells.Item(25, 17).Value = 22085

# Row 26
$ws.Cells.Item(26, 7).Value = 0.96
$ws.Cells.Item(26, 8).Value = 0.32
$ws.Cells.Item(26, 9).Value = 0.18
$ws.Cells.Item(26, 10).Value = 0.68
$ws.Cells.Item(26, 11).Value = 0.96
$ws.Cells.Item(26, 12).Value = 1.01
$ws.Cells.Item(26, 13).Value = 14.76
$ws.Cells.Item(26, 14).Value = 13.1
$ws.Cells.Item(26, 15).Value = 27.86
$ws.Cells.Item(26, 16).Value = 80
$ws.Cells.Item(26, 17).Value = 18082

# Row 27
$ws.Cells.Item(27, 7).Value = 0.6899999999999999
$ws.Cells.Item(27, 8).Value = 0.04
$ws.Cells.Item(27, 9).Value = 0.24
$ws.Cells.Item(27, 10).Value = 0.66
$ws.Cells.Item(27, 11).Value = 0.28
$ws.Cells.Item(27, 12).Value = 0.42
$ws.Cells.Item(27, 13).Value = 11.22
$ws.Cells.Item(27, 14).Value = 5.659999999999999
$ws.Cells.Item(27, 15).Value = 16.88
$ws.Cells.Item(27, 16).Value = 1746
$ws.Cells.Item(27, 17).Value = 21703

# Row 28
$ws.Cells.Item(28, 7).Value = 1.34
$ws.Cells.Item(28, 8).Value = 0.19
$ws.Cells.Item(28, 9).Value = 0.14
$ws.Cells.Item(28, 10).Value = 1.05
$ws.Cells.Item(28, 11).Value = 0.72
$ws.Cells.Item(28, 12).Value = 1.09
$ws.Cells.Item(28, 13).Value = 15.85
$ws.Cells.Item(28, 14).Value = 9.470000000000001
$ws.Cells.Item(28, 15).Value = 25.32
$ws.Cells.Item(28, 16).Value = 177
$ws.Cells.Item(28, 17).Value = 20068

# Row 29
$ws.Cells.Item(29, 7).Value = 0.74
$ws.Cells.Item(29, 8).Value = 0.29
$ws.Cells.Item(29, 9).Value = 0.12
$ws.Cells.Item(29, 10).Value = 0.75
$ws.Cells.Item(29, 11).Value = 0.7
$ws.Cells.Item(29, 12).Value = 0.87
$ws.Cells.Item(29, 13).Value = 13.47
$ws.Cells.Item(29, 14).Value = 8.85
$ws.Cells.Item(29, 15).Value = 22.32
$ws.Cells.Item(29, 16).Value = 144
$ws.Cells.Item(29, 17).Value = 20920

# Row 30
$ws.Cells.Item(30, 7).Value = 0.86
$ws.Cells.Item(30, 8).Value = 0.3
$ws.Cells.Item(30, 9).Value = 0.13
$ws.Cells.Item(30, 10).Value = 0.62
$ws.Cells.Item(30, 11).Value = 0.59
$ws.Cells.Item(30, 12).Value = 0.67
$ws.Cells.Item(30, 13).Value = 11.82
$ws.Cells.Item(30, 14).Value = 8.300000000000001
$ws.Cells.Item(30, 15).Value = 20.12
$ws.Cells.Item(30, 16).Value = 246
$ws.Cells.Item(30, 17).Value = 23086

# Row 31
$ws.Cells.Item(31, 7).Value = 0.61
$ws.Cells.Item(31, 8).Value = 0.28
$ws.Cells.Item(31, 9).Value = 0.18
$ws.Cells.Item(31, 10).Value = 0.9399999999999999
$ws.Cells.Item(31, 11).Value = 0.65
$ws.Cells.Item(31, 12).Value = 0.49
$ws.Cells.Item(31, 13).Value = 12.77
$ws.Cells.Item(31, 14).Value = 7.989999999999999
$ws.Cells.Item(31, 15).Value = 20.75999999999999
$ws.Cells.Item(31, 16).Value = 139
$ws.Cells.Item(31, 17).Value = 19561
